$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4305.952
$ws.Cells.Item(64, 9).Value = 3416.1333
$ws.Cells.Item(64, 10).Value = 6530.5
$ws.Cells.Item(64, 11).Value = 3416.1333
$ws.Cells.Item(64, 12).Value = 6530.5
$ws.Cells.Item(64, 13).Value = -3168.1333
$ws.Cells.Item(64, 14).Value = -7026.5

$ws.Cells.Item(67, 8).Value = 4305.952
$ws.Cells.Item(67, 9).Value = 3416.1333
$ws.Cells.Item(67, 10).Value = 6530.5
$ws.Cells.Item(67, 11).Value = 3416.1333
$ws.Cells.Item(67, 12).Value = 6530.5
$ws.Cells.Item(67, 13).Value = -2558.1333
$ws.Cells.Item(67, 14).Value = -8246.5

$ws.Cells.Item(98, 8).Value = 311850.22
$ws.Cells.Item(98, 9).Value = 386215.6
$ws.Cells.Item(98, 10).Value = 3765.1428
$ws.Cells.Item(98, 11).Value = 386215.6
$ws.Cells.Item(98, 12).Value = 3765.1428
$ws.Cells.Item(98, 13).Value = -384717.6
$ws.Cells.Item(98, 14).Value = -6761.1428

$ws.Cells.Item(113, 8).Value = 113021.555
$ws.Cells.Item(113, 9).Value = 168534.17
$ws.Cells.Item(113, 10).Value = 1996.3334
$ws.Cells.Item(113, 11).Value = 168534.17
$ws.Cells.Item(113, 12).Value = 1996.3334
$ws.Cells.Item(113, 13).Value = -165280.17
$ws.Cells.Item(113, 14).Value = -8504.3334

$ws.Cells.Item(122, 8).Value = 311850.22
$ws.Cells.Item(122, 9).Value = 386215.6
$ws.Cells.Item(122, 10).Value = 3765.1428
$ws.Cells.Item(122, 11).Value = 1158646.8
$ws.Cells.Item(122, 12).Value = 11295.4284
$ws.Cells.Item(122, 13).Value = -1156196.8
$ws.Cells.Item(122, 14).Value = -16195.4284

$ws.Cells.Item(132, 8).Value = 223471.75
$ws.Cells.Item(132, 9).Value = 259069.34
$ws.Cells.Item(132, 10).Value = 37573.223
$ws.Cells.Item(132, 11).Value = 777208.02
$ws.Cells.Item(132, 12).Value = 112719.669
$ws.Cells.Item(132, 13).Value = -774678.02
$ws.Cells.Item(132, 14).Value = -117779.669

$ws.Cells.Item(137, 8).Value = 1535.9584
$ws.Cells.Item(137, 9).Value = 886.1429000000001
$ws.Cells.Item(137, 10).Value = 1803.5294
$ws.Cells.Item(137, 11).Value = 2658.4287
$ws.Cells.Item(137, 12).Value = 5410.5882
$ws.Cells.Item(137, 13).Value = -108.4287000000004
$ws.Cells.Item(137, 14).Value = -10510.5882

$ws.Cells.Item(138, 8).Value = 5557652
$ws.Cells.Item(138, 9).Value = 1915.2122
$ws.Cells.Item(138, 10).Value = 8774131
$ws.Cells.Item(138, 11).Value = 5745.6366
$ws.Cells.Item(138, 12).Value = 26322393
$ws.Cells.Item(138, 13).Value = -605.6365999999998
$ws.Cells.Item(138, 14).Value = -26332673

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 11261.833
$ws.Cells.Item(31, 9).Value = 2267.75
$ws.Cells.Item(31, 10).Value = 29250
$ws.Cells.Item(31, 11).Value = 2267.75
$ws.Cells.Item(31, 12).Value = 29250
$ws.Cells.Item(31, 13).Value = -1973.75
$ws.Cells.Item(31, 14).Value = -29838

$ws.Cells.Item(32, 8).Value = 2096.3457
$ws.Cells.Item(32, 9).Value = 1247.2
$ws.Cells.Item(32, 11).Value = 1247.2
$ws.Cells.Item(32, 13).Value = -960.2

$ws.Cells.Item(61, 8).Value = 2086.743
$ws.Cells.Item(61, 9).Value = 1434.2759
$ws.Cells.Item(61, 10).Value = 5240.3335
$ws.Cells.Item(61, 11).Value = 1434.2759
$ws.Cells.Item(61, 12).Value = 5240.3335
$ws.Cells.Item(61, 13).Value = -1222.2759
$ws.Cells.Item(61, 14).Value = -5664.3335

$ws.Cells.Item(133, 8).Value = 49936.184
$ws.Cells.Item(133, 10).Value = 49936.184
$ws.Cells.Item(133, 12).Value = 49936.184
$ws.Cells.Item(133, 14).Value = -54996.184

$ws.Cells.Item(136, 8).Value = 2086.743
$ws.Cells.Item(136, 9).Value = 1434.2759
$ws.Cells.Item(136, 10).Value = 5240.3335
$ws.Cells.Item(136, 11).Value = 4302.8277
$ws.Cells.Item(136, 12).Value = 15721.0005
$ws.Cells.Item(136, 13).Value = -1752.8277
$ws.Cells.Item(136, 14).Value = -20821.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(97, 8).Value = 146356.75
$ws.Cells.Item(97, 9).Value = 146356.75
$ws.Cells.Item(97, 11).Value = 146356.75
$ws.Cells.Item(97, 13).Value = -145365.75

$ws.Cells.Item(99, 8).Value = 6212.375
$ws.Cells.Item(99, 9).Value = 1949.8334
$ws.Cells.Item(99, 11).Value = 1949.8334
$ws.Cells.Item(99, 13).Value = -451.8334

$ws.Cells.Item(134, 8).Value = 2635.1892
$ws.Cells.Item(134, 9).Value = 1622.7931
$ws.Cells.Item(134, 10).Value = 6305.125
$ws.Cells.Item(134, 11).Value = 4868.379300000001
$ws.Cells.Item(134, 12).Value = 18915.375
$ws.Cells.Item(134, 13).Value = -2333.379300000001
$ws.Cells.Item(134, 14).Value = -23985.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1297.0454
$ws.Cells.Item(31, 9).Value = 998.4
$ws.Cells.Item(31, 10).Value = 1545.9166
$ws.Cells.Item(31, 11).Value = 998.4
$ws.Cells.Item(31, 12).Value = 1545.9166
$ws.Cells.Item(31, 13).Value = -703.4
$ws.Cells.Item(31, 14).Value = -2135.9166

$ws.Cells.Item(34, 8).Value = 1297.0454
$ws.Cells.Item(34, 9).Value = 998.4
$ws.Cells.Item(34, 10).Value = 1545.9166
$ws.Cells.Item(34, 11).Value = 998.4
$ws.Cells.Item(34, 12).Value = 1545.9166
$ws.Cells.Item(34, 13).Value = -796.4
$ws.Cells.Item(34, 14).Value = -1949.9166

$ws.Cells.Item(58, 8).Value = 1836.174
$ws.Cells.Item(58, 9).Value = 583.5
$ws.Cells.Item(58, 10).Value = 3202.7273
$ws.Cells.Item(58, 11).Value = 583.5
$ws.Cells.Item(58, 12).Value = 3202.7273
$ws.Cells.Item(58, 13).Value = -380.5
$ws.Cells.Item(58, 14).Value = -3608.7273

$ws.Cells.Item(62, 8).Value = 17493.334
$ws.Cells.Item(62, 9).Value = 19366.666
$ws.Cells.Item(62, 10).Value = 10000
$ws.Cells.Item(62, 11).Value = 19366.666
$ws.Cells.Item(62, 12).Value = 10000
$ws.Cells.Item(62, 13).Value = -18742.666
$ws.Cells.Item(62, 14).Value = -11248

$ws.Cells.Item(65, 8).Value = 17493.334
$ws.Cells.Item(65, 9).Value = 19366.666
$ws.Cells.Item(65, 10).Value = 10000
$ws.Cells.Item(65, 11).Value = 96833.33
$ws.Cells.Item(65, 12).Value = 50000
$ws.Cells.Item(65, 13).Value = -93713.33
$ws.Cells.Item(65, 14).Value = -56240

$ws.Cells.Item(136, 8).Value = 1836.174
$ws.Cells.Item(136, 9).Value = 583.5
$ws.Cells.Item(136, 10).Value = 3202.7273
$ws.Cells.Item(136, 11).Value = 1750.5
$ws.Cells.Item(136, 12).Value = 9608.1819
$ws.Cells.Item(136, 13).Value = 799.5
$ws.Cells.Item(136, 14).Value = -14708.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 1915
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 1915
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 5745
$ws.Cells.Item(75, 14).Value = -7741
$ws.Cells.Item(75, 13).ClearContents()

$ws.Cells.Item(78, 8).Value = 1915
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 1915
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 17235
$ws.Cells.Item(78, 14).Value = -27219
$ws.Cells.Item(78, 13).ClearContents()

$ws.Cells.Item(113, 8).Value = 15151994
$ws.Cells.Item(113, 9).Value = 444
$ws.Cells.Item(113, 11).Value = 1332
$ws.Cells.Item(113, 13).Value = 838

$ws.Cells.Item(117, 8).Value = 1066.6666
$ws.Cells.Item(117, 10).Value = 2400
$ws.Cells.Item(117, 12).Value = 7200
$ws.Cells.Item(117, 14).Value = -14084

$ws.Cells.Item(121, 8).Value = 196.85715
$ws.Cells.Item(121, 9).Value = 179.66667
$ws.Cells.Item(121, 10).Value = 300
$ws.Cells.Item(121, 11).Value = 539.00001
$ws.Cells.Item(121, 12).Value = 900
$ws.Cells.Item(121, 13).Value = 770.99999
$ws.Cells.Item(121, 14).Value = -3520

$ws.Cells.Item(122, 8).Value = 667.2105
$ws.Cells.Item(122, 9).Value = 524.1429000000001
$ws.Cells.Item(122, 10).Value = 843.94116
$ws.Cells.Item(122, 11).Value = 4717.2861
$ws.Cells.Item(122, 12).Value = 7595.47044
$ws.Cells.Item(122, 13).Value = -2267.2861
$ws.Cells.Item(122, 14).Value = -12495.47044

$ws.Cells.Item(124, 8).Value = 6030
$ws.Cells.Item(124, 9).Value = 6030
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 11).Value = 18090
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 13).Value = -13180
$ws.Cells.Item(124, 14).ClearContents()

$ws.Cells.Item(125, 8).Value = 2751.0833
$ws.Cells.Item(125, 9).Value = 2247.5
$ws.Cells.Item(125, 10).Value = 3002.875
$ws.Cells.Item(125, 11).Value = 6742.5
$ws.Cells.Item(125, 12).Value = 9008.625
$ws.Cells.Item(125, 13).Value = -1822.5
$ws.Cells.Item(125, 14).Value = -18848.625

$ws.Cells.Item(131, 8).Value = 2543.6582
$ws.Cells.Item(131, 10).Value = 2768.8872
$ws.Cells.Item(131, 12).Value = 8306.661599999999
$ws.Cells.Item(131, 14).Value = -18386.6616

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 855779.75
$ws.Cells.Item(122, 9).Value = 1111903.8
$ws.Cells.Item(122, 10).Value = 2033.3334
$ws.Cells.Item(122, 11).Value = 3335711.4
$ws.Cells.Item(122, 12).Value = 6100.0002
$ws.Cells.Item(122, 13).Value = -3333261.4
$ws.Cells.Item(122, 14).Value = -11000.0002

$ws.Cells.Item(138, 8).Value = 64100
$ws.Cells.Item(138, 10).Value = 64100
$ws.Cells.Item(138, 12).Value = 64100
$ws.Cells.Item(138, 14).Value = -74380

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3515.3845
$ws.Cells.Item(7, 9).Value = 3200
$ws.Cells.Item(7, 10).Value = 3541.6667
$ws.Cells.Item(7, 11).Value = 3200
$ws.Cells.Item(7, 12).Value = 3541.6667
$ws.Cells.Item(7, 13).Value = -3088
$ws.Cells.Item(7, 14).Value = -3765.6667

$ws.Cells.Item(126, 8).Value = 3515.3845
$ws.Cells.Item(126, 9).Value = 3200
$ws.Cells.Item(126, 10).Value = 3541.6667
$ws.Cells.Item(126, 11).Value = 9600
$ws.Cells.Item(126, 12).Value = 10625.0001
$ws.Cells.Item(126, 13).Value = -7130
$ws.Cells.Item(126, 14).Value = -15565.0001

$ws.Cells.Item(132, 8).Value = 4074.875
$ws.Cells.Item(132, 9).Value = 3399
$ws.Cells.Item(132, 11).Value = 10197
$ws.Cells.Item(132, 13).Value = -7667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(104, 8).Value = 44500
$ws.Cells.Item(104, 10).Value = 44500
$ws.Cells.Item(104, 12).Value = 44500
$ws.Cells.Item(104, 14).Value = -51488

$ws.Cells.Item(132, 8).Value = 13159868
$ws.Cells.Item(132, 9).Value = 17858744
$ws.Cells.Item(132, 10).Value = 3014.4
$ws.Cells.Item(132, 11).Value = 53576232
$ws.Cells.Item(132, 12).Value = 9043.200000000001
$ws.Cells.Item(132, 13).Value = -53573702
$ws.Cells.Item(132, 14).Value = -14103.2

$ws.Cells.Item(136, 8).Value = 9834472
$ws.Cells.Item(136, 9).Value = 12859542
$ws.Cells.Item(136, 11).Value = 38578626
$ws.Cells.Item(136, 13).Value = -38576076
